$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Mensaje" column (B) and its sample messages are gone - clear header + data.
$ws.Range("B1:B3").ClearContents()

# Row 2 now holds the phone number itself (with country code) as text, not a number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "+525529282277"

# Row 3 loses its old phone number; the cell remains but picks up an underline style
# (leftover formatting from the template), so re-create it as an empty, underlined cell.
$ws.Range("A3").ClearContents()
$ws.Range("A3").Font.Underline = $true

# Two further blank rows are present below (kept as real, empty, styled cells so the
# sheet's used range grows to A1:A5, matching the saved template).
$ws.Range("A1").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths from the template (A narrow for phone numbers, B wide for messages).
$ws.Range("A1").ColumnWidth = 12.25
$ws.Range("B1").ColumnWidth = 84.75

# The saved view now has A3 selected instead of B4.
$ws.Range("A3").Select() | Out-Null
